$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row with two new columns (P1, Q1), matching the style of
# the existing header cells (bold/centered/bordered style used by O1).
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$hdr = $ws.Range("P1:Q1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Update data rows 2-25: swap values in columns I, K, M, O and add new
# columns P, Q (each filled with 2).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2
}
